$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.185.45"
$ws.Range("E2").Value = "  -0.70%  "

$ws.Range("D3").Value = "2.450.16"
$ws.Range("E3").Value = "  +0.28%  "

$ws.Range("E4").Value = "  -0.09%  "

$ws.Range("D5").Value = "'582.92"
$ws.Range("E5").Value = "  +1.27%  "

$ws.Range("D6").Value = "'143.22"
$ws.Range("E6").Value = "  -0.83%  "

$ws.Range("E7").Value = "  +0.04%  "

$ws.Range("E8").Value = "  +0.30%  "

$ws.Range("D9").Value = "2.446.13"
$ws.Range("E9").Value = "  +0.27%  "

$ws.Range("E10").Value = "  +2.38%  "

$ws.Range("E11").Value = "  +2.40%  "

$ws.Range("E12").Value = "  -0.01%  "

$ws.Range("E13").Value = "  -2.64%  "

$ws.Range("D14").Value = "'26.46"
$ws.Range("E14").Value = "  -0.98%  "

$ws.Range("E15").Value = "  +0.32%  "

$ws.Range("D16").Value = "2.886.06"
$ws.Range("E16").Value = "  +0.12%  "

$ws.Range("D17").Value = "62.098.81"
$ws.Range("E17").Value = "  -0.55%  "

$ws.Range("D18").Value = "2.443.06"
$ws.Range("E18").Value = "  +0.02%  "

$ws.Range("D19").Value = "'10.69"
$ws.Range("E19").Value = "  -4.30%  "

$ws.Range("E20").Value = "  +1.07%  "

$ws.Range("D21").Value = "'326.77"
$ws.Range("E21").Value = "  -0.63%  "

$ws.Range("E22").Value = "  -1.29%  "

$ws.Range("E23").Value = "  +0.14%  "

$ws.Range("D24").Value = "'1.91"
$ws.Range("E24").Value = "  -5.90%  "

$ws.Range("D25").Value = "'65.68"
$ws.Range("E25").Value = "  +0.15%  "

$ws.Range("D26").Value = "'9.13"
$ws.Range("E26").Value = "  -0.95%  "

$ws.Range("D27").Value = "'599.65"
$ws.Range("E27").Value = "  -5.44%  "

$ws.Range("D28").Value = "0.0₃0971"
$ws.Range("E28").Value = "  +0.53%  "

$ws.Range("E29").Value = "  +0.16%  "

$ws.Range("E30").Value = "  -0.28%  "

$ws.Range("D31").Value = "'1.42"
$ws.Range("E31").Value = "  -1.16%  "

$ws.Range("D32").Value = "'8.00"
$ws.Range("E32").Value = "  -1.16%  "

$ws.Range("D33").Value = "'1.90"
$ws.Range("E33").Value = "  +0.82%  "

$ws.Range("D34").Value = "'0.135"
$ws.Range("E34").Value = "  -1.15%  "

$ws.Range("D35").Value = "'4.89"
$ws.Range("E35").Value = "  -2.44%  "

$ws.Range("E36").Value = "  +0.15%  "

$ws.Range("E37").Value = "  -2.18%  "

$ws.Range("E39").Value = "  +4.52%  "

$ws.Range("B40").Value = "RenderToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D40").Value = "'5.31"
$ws.Range("E40").Value = "  +0.57%  "

$ws.Range("B41").Value = "EthereumClassic"
$ws.Range("C41").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D41").Value = "'5.31"
$ws.Range("E41").Value = "  -0.62%  "

$ws.Range("D42").Value = "'43.19"
$ws.Range("E42").Value = "  +2.13%  "

$ws.Range("E43").Value = "  -1.58%  "

$ws.Range("D45").Value = "'2.53"
$ws.Range("E45").Value = "  +0.81%  "

$ws.Range("D46").Value = "0.0₆0279"
$ws.Range("E46").Value = "  +24.17%  "

$ws.Range("D47").Value = "'141.99"
$ws.Range("E47").Value = "  -2.62%  "

$ws.Range("E48").Value = "  -2.71%  "

$ws.Range("D49").Value = "'0.602"
$ws.Range("E49").Value = "  +0.61%  "

$ws.Range("E50").Value = "  -1.18%  "

$ws.Range("D51").Value = "'19.89"
$ws.Range("E51").Value = "  +0.81%  "
